{"js": "// The edit removes the trailing \"footer\" block that was appended to the\n// document body: an empty paragraph, the \"Ver no Jupiter...\" paragraph and\n// the \"\u00a9 2020 ...\" paragraph. They sit right after the paragraph that\n// contains \"LOM3094: Processamento de Materiais Met\u00e1licos II (Requisito)\"\n// and right before the pre-existing trailing empty paragraph / page-break\n// paragraph, which must be kept untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the two textual paragraphs that must be deleted.\nlet jupiterIndex = -1;\nlet copyrightIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text;\n  if (jupiterIndex === -1 && t.indexOf(\"Ver no Jupiter\") !== -1) {\n    jupiterIndex = i;\n  }\n  if (copyrightIndex === -1 && t.indexOf(\"Powered by Jekyll\") !== -1) {\n    copyrightIndex = i;\n  }\n}\n\nconst toDelete = [];\nif (jupiterIndex !== -1) {\n  // The empty paragraph immediately preceding \"Ver no Jupiter...\" is part\n  // of the same trailing block and must go too.\n  const prev = items[jupiterIndex - 1];\n  if (prev && prev.text.trim() === \"\") {\n    toDelete.push(prev);\n  }\n  toDelete.push(items[jupiterIndex]);\n}\nif (copyrightIndex !== -1) {\n  toDelete.push(items[copyrightIndex]);\n}\n\n// Delete from the end backwards so earlier indices stay valid.\nfor (let i = toDelete.length - 1; i >= 0; i--) {\n  toDelete[i].delete();\n}\nawait context.sync();\n", "ps1": "# The edit removes the trailing \"footer\" block that was appended to the\n# document body: an empty paragraph, the \"Ver no Jupiter...\" paragraph and\n# the \"(c) 2020 ...\" paragraph. They sit right after the paragraph that\n# contains \"LOM3094: Processamento de Materiais Metalicos II (Requisito)\"\n# and right before the pre-existing trailing empty paragraph / page-break\n# paragraph, which must be kept untouched.\n\n$d = $word.ActiveDocument\n\n$jupiterIndex = -1\n$copyrightIndex = -1\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($jupiterIndex -eq -1 -and $t -like \"*Ver no Jupiter*\") {\n        $jupiterIndex = $i\n    }\n    if ($copyrightIndex -eq -1 -and $t -like \"*Powered by Jekyll*\") {\n        $copyrightIndex = $i\n    }\n}\n\n$toDelete = @()\n\nif ($jupiterIndex -ne -1) {\n    $prevText = $d.Paragraphs.Item($jupiterIndex - 1).Range.Text.Trim()\n    if ($prevText -eq \"\") {\n        $toDelete += ($jupiterIndex - 1)\n    }\n    $toDelete += $jupiterIndex\n}\nif ($copyrightIndex -ne -1) {\n    $toDelete += $copyrightIndex\n}\n\n# Delete highest index first so earlier indices stay valid.\n$sorted = $toDelete | Sort-Object -Descending\nforeach ($idx in $sorted) {\n    $d.Paragraphs.Item($idx).Range.Delete()\n}\n"}
